$wb = $excel.ActiveWorkbook

# --- Fill in row 5 of the "ProviderOptions" sheet with the new test-case data ---
$ws = $wb.Worksheets.Item("ProviderOptions")

$ws.Range("A5").Value = "testT4275"
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 4
$ws.Range("D5").Value = "Click"
$ws.Range("E5").Value = "Click"
$ws.Range("F5").Value = "Click"
$ws.Range("G5").Value = "n/a"
$ws.Range("H5").Value = "n/a"
$ws.Range("I5").Value = "n/a"
$ws.Range("J5").Value = "No"
$ws.Range("K5").Value = "n/a"
$ws.Range("L5").Value = "n/a"
$ws.Range("M5").Value = "n/a"
$ws.Range("N5").Value = "n/a"
$ws.Range("O5").Value = "Click"
$ws.Range("P5").Value = "ProviderOption_ID1"
$ws.Range("Q5").Value = "n/a"
$ws.Range("R5").Value = "n/a"
$ws.Range("S5").Value = "n/a"
$ws.Range("T5").Value = "n/a"

# U5 is a brand-new cell (row 5 previously stopped at column T), so give it the
# same style as the rest of the row (copy formatting from T5) before setting its value.
$ws.Range("T5").Copy()
$ws.Range("U5").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("U5").Value = "Placement Decision"

# --- Activate "ProviderOptions" and select C5 (moves tabSelected + selection here) ---
$ws.Activate()
$ws.Range("C5").Select()
